$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 held the phone number 791000072 as a real number; it now becomes a
# text value (quote-prefixed) with the same "791000072" content.
$ws.Range("A2").Value = "'791000072"

# C2/D2 held the password "Kalemon12345@" (also used as the hyperlink
# display text); update it to the new password text.
$ws.Range("C2").Value = "Kalemon12345678@"
$ws.Range("D2").Value = "Kalemon12345678@"

# Update the saved cursor/selection position shown in the sheet view.
$ws.Range("D6").Select() | Out-Null
